$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (A1:D8) needs to be sorted in ascending order by column A
# (time), which reorders the data rows (A2:D8) while keeping the header
# row (row 1) fixed.
$dataRange = $ws.Range("A1:D8")
$keyRange = $ws.Range("A1")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 1
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.SortMethod = 1
$ws.Sort.Apply()

$wb.Save()
